# Regenerate the "K" column (G) of the save_data sheet.
# The upstream data pipeline recomputed the strike-count ("K", formerly
# "Strike#") values together with std/mean statistics and wrote the new
# s_vals into column G for every data row (rows 2-64).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(2, 1, 1, 1, 2, 1, 0, 0, 0, 1, 0, 2, 0, 1, 1, 0, 3, 1, 0, 1, 3, 1, 1, 0, 1, 1, 0, 0, 4, 2, 0, 2, 0, 0, 1, 2, 1, 0, 1, 0, 0, 1, 0, 0, 1, 0, 0, 1, 1, 1, 0, 0, 2, 0, 1, 0, 2, 1, 2, 2, 1, 2, 0)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
